# "recording the war data"
# Insert a new scene row (row 11) above the existing "villageScene1" entry
# (which was row 11 and is now pushed down to row 12, together with every
# row below it). The new row re-uses the same scene/asset data as
# villageScene1 but records a new Id ("0") and a much larger
# MaxGroupPlayers value (500000) — the "war" scene.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 11:16 down to 12:17, leaving a blank row 11 behind.
$ws.Rows("11:11").Insert()

# The freshly inserted row inherits stray formatting from the row above
# (the header band) — strip that back to the plain, unstyled look the
# data rows below it use.
$ws.Range("A11:Q11").ClearFormats()

# Columns that hold text (including numeric-looking text like "0" or
# "45,180") need an explicit Text number format so COM/Excel doesn't
# silently reinterpret them as numbers.
$ws.Range("A11").NumberFormat = "@"
$ws.Range("F11").NumberFormat = "@"
$ws.Range("G11").NumberFormat = "@"
$ws.Range("I11").NumberFormat = "@"
$ws.Range("J11").NumberFormat = "@"
$ws.Range("N11").NumberFormat = "@"
$ws.Range("O11").NumberFormat = "@"
$ws.Range("P11").NumberFormat = "@"
$ws.Range("Q11").NumberFormat = "@"

$ws.Range("A11").Value = "0"
$ws.Range("B11").Value = "villageScene1"
$ws.Range("C11").Value = "villageScene1"
$ws.Range("D11").Value = 100
$ws.Range("E11").Value = 500000
$ws.Range("F11").Value = "../NFDataCfg/Ini/Scene/1.xml"
$ws.Range("G11").Value = "89,104,0"
$ws.Range("H11").Value = 500
$ws.Range("I11").Value = "Sources/Music/Town"
$ws.Range("J11").Value = "Sources/Music/Town"
$ws.Range("K11").Value = 0
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = 0
$ws.Range("N11").Value = "UI/ChronoBlade_forest_wallpaper"
$ws.Range("O11").Value = "0,8,7"
$ws.Range("P11").Value = "45,180"
$ws.Range("Q11").Value = "../NFDataCfg/Ini/Navigation/srv_demo.navmesh"

# Match the author's final cursor position noted in the saved workbook.
$ws.Range("L10").Select() | Out-Null
